$wb = $excel.ActiveWorkbook

$wsAdd = $wb.Worksheets.Item("Add Devices")
$wsDel = $wb.Worksheets.Item("Delete Devices")

# Update the boolean cells (previously text "NA") on the "Add Devices" sheet
$wsAdd.Range("L10:M13").Value = $false

# Update the boolean cells (previously text "NA") on the "Delete Devices" sheet
$wsDel.Range("L10:M10").Value = $false

# Update the selection on "Delete Devices" first (it is currently the active
# sheet), then switch the active sheet to "Add Devices" and set its
# selection last so it ends up being the active tab/window.
[void]$wsDel.Range("L10").Select()

[void]$wsAdd.Activate()
[void]$wsAdd.Range("L10:M13").Select()
